# Apply "ran scans. Added words. Slight calculation change" edit:
#  - Fill in the previously-blank X5/Y5 cells on row 5.
#  - Append a brand-new row 6 with a fresh scan result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- complete row 5 (previously only went through column W) ---
$ws.Range("X5").Value = 0.6499990000000011
$ws.Range("Y5").Value = "Up"

# --- new row 6 ---
# Copy the date / percentage number formats down from row 5 first (so we
# reuse the existing style records instead of minting new numFmt entries),
# then fill in the values.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("S5:T5").Copy()
$ws.Range("S6:T6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A6").Value = 42647.884305555555
$ws.Range("B6").Value = 11
$ws.Range("C6").Value = "Buy"
$ws.Range("D6").Value = 34
$ws.Range("E6").Value = 15769
$ws.Range("F6").Value = 801
$ws.Range("G6").Value = 60
$ws.Range("H6").Value = 38
$ws.Range("I6").Value = 80
$ws.Range("J6").Value = 19
$ws.Range("K6").Value = 24355
$ws.Range("L6").Value = 202
$ws.Range("M6").Value = 127
$ws.Range("N6").Value = 67
$ws.Range("O6").Value = 16
$ws.Range("P6").Value = "Named"
$ws.Range("Q6").Value = 46.357611069683557
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0.0591
$ws.Range("T6").Value = -0.042099999999999999
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = "N/A"
$ws.Range("W6").Value = 0
